$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Fix the typo'd product/test-case name in B1 of both sheets
# (old text had a stray space before "Repayment")
$newName = "1015-MS-EI-DB-SAR-REC-NON-RNI-CTRFD-SAR-MD-TR-1-LateRepayment"
$wsInput.Range("B1").Value  = $newName
$wsOutput.Range("B1").Value = $newName

# Re-apply the formatting those header cells carry elsewhere on the sheet
# (bold-green header look used throughout column B)
$wsInput.Range("B1").Font.Name    = "Arial"
$wsInput.Range("B1").Font.Size    = 10
$wsInput.Range("B1").Interior.Color = $wsInput.Range("B23").Interior.Color

$wsOutput.Range("B1").Font.Name    = "Arial"
$wsOutput.Range("B1").Font.Size    = 10
$wsOutput.Range("B1").Interior.Color = $wsInput.Range("B23").Interior.Color

# Make ProductLoanInput the active/selected sheet again, scrolled back to the
# top with the cursor on B1 (previously ProductLoanOutput was the active tab
# and ProductLoanInput was scrolled down to A27). Set ProductLoanOutput's
# selection to B1 first so activating ProductLoanInput afterwards is what
# sticks as the workbook's active tab.
$wsOutput.Range("B1").Select() | Out-Null

$wsInput.Activate() | Out-Null
$wsInput.Range("B1").Select() | Out-Null
